$p = $ppt.ActivePresentation

# --- Slide 13: remove the "controleren van de configuratie ..." bullet paragraph ---
$slide13 = $p.Slides.Item(13)
$shape13 = $slide13.Shapes.Item(2)
$tr13 = $shape13.TextFrame.TextRange
$fullText13 = $tr13.Text
$target13 = "controleren van de configuratie op aanwezigheid van bekende kwetsbaarheden,`r"
$idx13 = $fullText13.IndexOf($target13)
if ($idx13 -ge 0) {
    $sub13 = $tr13.Characters($idx13 + 1, $target13.Length)
    $sub13.Delete()
}

# --- Slide 19: update M23 title to mention "kennis van en" ---
$slide19 = $p.Slides.Item(19)
$titleShape = $slide19.Shapes.Item(1)
$titleTr = $titleShape.TextFrame.TextRange
$titleFull = $titleTr.Text
$titleAll = $titleTr.Characters(1, $titleFull.Length)
$titleAll.Text = "M23: Het project zorgt voor de aanwezigheid van kennis van en ervaring met de Kwaliteitsaanpak"

# --- Slide 19: extend body paragraph with the extra explanation sentence ---
$bodyShape = $slide19.Shapes.Item(2)
$bodyTr = $bodyShape.TextFrame.TextRange
$bodyFull = $bodyTr.Text
$bodyAll = $bodyTr.Characters(1, $bodyFull.Length)
$bodyAll.Text = "De software delivery manager zorgt ervoor dat bij nieuwe projecten wordt gestart met ten minste twee projectleden die bekend zijn met de Kwaliteitsaanpak. Projectleden die nog niet bekend zijn met de Kwaliteitsaanpak krijgen uitleg over de inhoud en achtergrond van de Kwaliteitsaanpak."
